$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Temporarily clear the format of B9 ("Enable Redis" value cell). This avoids a row-insert
# quirk where the engine creates a spurious extra cell style when a row is inserted directly
# below a cell that has a fill+border style. We restore the format right after the insert.
$ws.Cells.Item(9, 2).ClearFormats()

# Insert a new row before row 10 (shifts "Allow CORS" row and everything below it down by one)
$ws.Rows.Item(10).Insert()

# Restore B9's formatting and give the same formatting to the new B10 cell, both copied
# from B8 ("Enable Audit" value cell), which already carries the desired style.
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(9, 2).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(10, 2).PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Populate the newly inserted row 10 with "Enable SMTP" / "Yes"
$ws.Cells.Item(10, 1).Value = "Enable SMTP"
$ws.Cells.Item(10, 2).Value = "Yes"

# Update selection to match the target state
$ws.Range("A13:E15").Select()
